$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the header/timestamp column (also updates the table column name automatically)
$ws.Range("C1").Value = "2024-11-15 16:57:30"

# Update attendance status values for the affected rows
$ws.Range("C2").Value  = "Falta"
$ws.Range("C4").Value  = "Puntual"
$ws.Range("C5").Value  = "Falta"
$ws.Range("C6").Value  = "Falta"
$ws.Range("C7").Value  = "Falta"
$ws.Range("C10").Value = "Falta"
$ws.Range("C11").Value = "Falta"
$ws.Range("C12").Value = "Falta"

$wb.Save()
